$wb = $excel.ActiveWorkbook

# The handback transform failed for the 4500359d-... file (row 3). Update
# the Status column everywhere it is reported (Overview + each locale
# sheet) to reflect the failure, and record the error detail (column K)
# on each locale sheet describing the mismatched file name.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("K3").Value = "Handback file name: z1qzfxwr.gn3 is different with handoff file name: 4500359d-6849-4cbe-9929-3bb9708bddde.510cfa9a4b4a92497952e7f6bb990a461b4a94ba.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("K3").Value = "Handback file name: z1qzfxwr.gn3 is different with handoff file name: 4500359d-6849-4cbe-9929-3bb9708bddde.510cfa9a4b4a92497952e7f6bb990a461b4a94ba.de-de."

$wb.Save()
